$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 495, shifting existing rows 495:568 down to 496:569.
$ws.Rows.Item(495).Insert()

# Populate the newly inserted row 495 with the new record. All columns
# except D (Fecha) and M (Volumen) mirror the row that was previously at 495
# (now shifted to 496): Perú / Primera / $/bandeja 4 kilos / etc.
$ws.Range("A495").Value = 10
$ws.Range("B495").Value = "Vega Modelo de Temuco"
$ws.Range("C495").Value = "La Araucanía"
$ws.Range("D495").Value = 45077
$ws.Range("E495").Value = 9
$ws.Range("F495").Value = "Fruta"
$ws.Range("G495").Value = 100108
$ws.Range("H495").Value = "Tropicales y subtropicales"
$ws.Range("I495").Value = 100108002
$ws.Range("J495").Value = "Mango"
$ws.Range("K495").Value = "Sin especificar"
$ws.Range("L495").Value = "Primera"
$ws.Range("M495").Value = 500
$ws.Range("N495").Value = 8000
$ws.Range("O495").Value = 8000
$ws.Range("P495").Value = 8000
$ws.Range("Q495").Value = "`$/bandeja 4 kilos"
$ws.Range("R495").Value = "Perú"
$ws.Range("S495").Value = 2000
$ws.Range("T495").Value = 4
